# Update report co so: refresh data in both sheets (DOANH SỐ CÁ NHÂN / CHI TIÊU)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet 1: DOANH SỐ CÁ NHÂN  (new data set, rows 2-16, header row 1 unchanged)
# Columns: A Name | B..J numeric figures
# ---------------------------------------------------------------------------
$sheet1Rows = @(
    @("Bác Sĩ Ngoài",            0,         0,          5500000,  0, 0,  0,      0, 0, 0),
    @("Bác Sĩ Thảo",             0,         0,          25000000, 0, 0,  0,      0, 0, 0),
    @("CTV Ngoài",               500000,    0,          0,        0, 0,  0,      0, 0, 500000),
    @("Kha Như Huỳnh ",          800000,    0,          800000,   0, 5,  200000, 0, 0, 0),
    @("Lâm Hoàng Phú",           0,         0,          0,        0, 4,  200000, 0, 0, 0),
    @("Lâm Thị Mỹ Hằng",         12500000,  0,          13500000, 0, 0,  0,      0, 0, 19000000),
    @("Lê Đình Hậu",             7000000,   0,          0,        0, 0,  0,      0, 0, 2000000),
    @("Nguyễn Hoàng Yến Quyên",  4500000,   0,          42200000, 0, 0,  0,      0, 0, 0),
    @("Nguyễn Phúc Nam",         3000000,   0,          0,        0, 0,  0,      0, 0, 800000),
    @("Phạm Thanh Hoàng",        0,         0,          17730000, 0, 0,  0,      0, 0, 0),
    @("Thạch Hoàng Nhân",        66000000,  0,          0,        0, 0,  0,      0, 0, 12000000),
    @("Đào Vương Anh",           0,         0,          0,        0, 1,  0,      0, 0, 0),
    @("Đặng Ngọc Mai",           0,         0,          4500000,  0, 0,  0,      0, 0, 0),
    @("Đỗ Thị Huyền Trân",       21430000,  6000000,    0,        0, 0,  0,      0, 0, 0),
    @("Tổng",                    115730000, 6000000,    109230000,0, 10, 400000, 0, 0, 34300000)
)

$r = 2
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    for ($c = 1; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: CHI TIÊU  (new data set, rows 2-12, header row 1 unchanged)
# Columns: A Category | B amount
# ---------------------------------------------------------------------------
$sheet2Rows = @(
    @("Chi Phí CTV",                    20730000),
    @("Chi Phí Hạ Tầng",                4544000),
    @("Chi Phí Sinh Hoạt Tại Cơ Sở",    3768000),
    @("Chi Phí Vận Hành",               4205000),
    @("Chí Phí Bác Sĩ Ngoài",           6250000),
    @("Phúc lợi công ty",               600000),
    @("Tiền Thuế",                      0),
    @("Trang thiết bị Y Tế",            3540000),
    @("Ứng Lương",                      6699000),
    @("Blank",                          5000000),
    @("Tổng cộng",                      55336000)
)

$r = 2
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
